# CPSC224 Space Race Project Plan - weekly update (4/6 edit)
# Work against the "WeekApr1" sheet (tab index 2 / WeekApr1), which is the
# active sheet in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeekApr1")

# --- Row 3: "Priliminay class design" task ---------------------------------
# Re-assign from Zach to the Group; mark complete with a completion date and
# actual effort logged.
$ws.Range("B3").Value = "Group"
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = "4/4/2018"
$ws.Range("G3").Value = 0.5

# --- Row 4: "Desigining the classes needed" task ----------------------------
# Mark complete with a completion date and actual effort logged.
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = "4/5/2018"
$ws.Range("G4").Value = 1

# --- Row 5: "Create UML diagram" task ---------------------------------------
# Re-assign from Andrew to Zach; mark complete with a completion date.
$ws.Range("B5").Value = "Zach"
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = "4/6/2018"

# --- Row 6: "UI Mockup" task -------------------------------------------------
# Mark complete with a completion date and actual effort logged.
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = "4/6/2018"
$ws.Range("G6").Value = 1.5

# --- Row 7: new task "Work on resource icons" -------------------------------
$ws.Range("A7").Value = "Work on resource icons"
$ws.Range("B7").Value = "Andrew"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "4/13/2018"
$ws.Range("E7").Value = 40

# Give F7 the same date-formatted look as the other Completion Date cells in
# this column, without actually filling in a completion date yet (task is
# still in progress). Copy formatting only from an already-formatted cell.
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Sheet-level bits --------------------------------------------------------
# Cursor was left on E8 when the workbook was last saved.
$ws.Range("E8").Select() | Out-Null

# Set the sheet up for portrait printing.
$ws.PageSetup.Orientation = 1
